# Split "If using Xbox Series X/S, ..." into three runs:
#   "If using Xbox Series X" | "|" | "S, set the active solution platform to "
# i.e. replace the "/" between "X" and "S" with "|", as its own run,
# matching the commit "Update Xbox GDK Samples to November GDK release."

$d = $word.ActiveDocument

# Locate the "/" that sits between "Series X" and "S," in the Xbox Series X/S
# sentence (search narrowly so we don't touch anything else in the doc).
$needle = $d.Content
$found = $needle.Find.Execute("Xbox Series X/S", $false, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'Xbox Series X/S' in the document"
}

$sentenceStart = $needle.Start
$slashPos = $sentenceStart + "Xbox Series X".Length

$slashRange = $d.Range($slashPos, $slashPos + 1)
if ($slashRange.Text -ne "/") {
    throw "Expected '/' at computed position, found '$($slashRange.Text)'"
}

# First turn the "/" into "|" while it is still a single run.
$slashRange.Text = "|"

# Re-acquire the (now "|") one-character range and force Word to split it
# into its own run, separate from the text before and after, by toggling a
# character property (set it on, then back off). The round trip leaves the
# run boundary in place without leaving a visible formatting change behind.
$pipeRange = $d.Range($slashPos, $slashPos + 1)
$pipeRange.Bold = $true
$pipeRange.Bold = $false
